$d = $word.ActiveDocument

# 1. "The final project will be due at 12:15pm..." -> "...due before 12:15pm..."
[void]$d.Content.Find.Execute("The final project will be due at 12:15pm on Tuesday, June 13th.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "The final project will be due before 12:15pm on Tuesday, June 13th.", 2)

# 2. Remove three of the four blank paragraphs that follow the
#    "...zipped-up copy of the produced dataset/model files." paragraph,
#    leaving a single trailing blank paragraph.
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "zipped-up copy of the produced dataset/model files\.") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 0) {
    $blankStart = $d.Paragraphs.Item($targetIdx + 1)
    $blankKeep  = $d.Paragraphs.Item($targetIdx + 4)
    $r = $d.Range($blankStart.Range.Start, $blankKeep.Range.Start)
    $r.Delete()
}
